$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph "Problema de transferencia de calor por conducción."
# -> "Problema de " + bold("transferencia de calor por conducción")
#    wrapped in bookmark _Hlk84178026 + "."
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(6)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)
$r1.Find.Execute("transferencia de calor por conducción")
$r1.Bold = $true
$d.Bookmarks.Add("_Hlk84178026", $r1)

# ------------------------------------------------------------------
# Paragraph "Cálculo de transferencia de calor por convección para:"
# -> "Cálculo de " + bold("transferencia de calor por convección")
#    wrapped in bookmark _Hlk84178179 + " para:"
# ------------------------------------------------------------------
$p2 = $d.Paragraphs(7)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$r2.Find.Execute("transferencia de calor por convección")
$r2.Bold = $true
$d.Bookmarks.Add("_Hlk84178179", $r2)

# ------------------------------------------------------------------
# Paragraph "Convección forzada" -> whole paragraph bold
# ------------------------------------------------------------------
$p3 = $d.Paragraphs(8)
$p3.Range.Bold = $true
$p3.Range.BoldBi = $true

# ------------------------------------------------------------------
# Paragraph "Convección natural" -> whole paragraph bold
# ------------------------------------------------------------------
$p4 = $d.Paragraphs(9)
$p4.Range.Bold = $true
$p4.Range.BoldBi = $true

Write-Output "edit applied"
